$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old content (rows 4:5, cols A:C) entirely.
$ws.Range("A4:C5").Clear()

# Write the new, smaller table starting at A1.
$ws.Range("A1").Value = "_id"
$ws.Range("B1").Value = "age"
$ws.Range("A2").Value = 140
$ws.Range("B2").Value = 22

$ws.Range("B2").Select()
